# "system architecture annd summary added"
#
# The slide previously held a flat set of shapes: two arrow connectors,
# two rotated callout textboxes, the "Diagram ..." title textbox, and
# three screenshot pictures. The edit removes the title textbox
# ("TextBox 9") entirely and collects the remaining seven shapes
# (the connectors, the two callout textboxes, and the three pictures)
# into a single new group shape ("Group 1").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Remove the standalone "Diagram textual_challenge_mult_choice_response_results"
# title textbox (id 10) - it is not part of the new group.
$title = $s.Shapes.Item("TextBox 9")
$title.Delete()

# Group the remaining diagram shapes together.
$names = @(
    "Straight Arrow Connector 26",
    "Straight Arrow Connector 27",
    "TextBox 29",
    "TextBox 31",
    "Picture 10",
    "Picture 11",
    "Picture 12"
)
$range = $s.Shapes.Range($names)
$group = $range.Group()
$group.Name = "Group 1"
